$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting all existing columns
# (and the data validation range) one position to the right.
$ws.Columns("A").Insert()

# The newly inserted column A has no formatting; copy the formatting
# from the (now adjacent) column B for the two used rows so the new
# cells pick up the same header / data styles as the rest of the row.
$ws.Range("B1:B2").Copy()
$ws.Range("A1:A2").PasteSpecial(-4122)  # xlPasteFormats

# Set the new column's width (Excel's ColumnWidth property is in
# character units and gets padded by ~5/6 of a character internally,
# so subtract that padding to land on the desired stored width of 23).
$ws.Columns("A").ColumnWidth = 23 - 5/6

# Populate the new index column.
$ws.Range("A1").Value = "INDEX (DO NOT MODIFY)"
$ws.Range("A2").Value = 1

# Upper-case all of the original header labels, which now live in
# columns B through X (column Y is the untouched "status as of" column).
$ws.Range("B1").Value = "REGION"
$ws.Range("C1").Value = "DIVISION"
$ws.Range("D1").Value = "SCHOOL ID"
$ws.Range("E1").Value = "SCHOOL NAME"
$ws.Range("F1").Value = "MUNICIPALITY"
$ws.Range("G1").Value = "LEG DISTRICT"
$ws.Range("H1").Value = "NO. OF SITES"
$ws.Range("I1").Value = "SCOPE OF WORK"
$ws.Range("J1").Value = "TOTAL ALLOCATION"
$ws.Range("K1").Value = "CONTRACT AMOUNT"
$ws.Range("L1").Value = "STATUS"
$ws.Range("M1").Value = "PERCENTAGE OF COMPLETION"
$ws.Range("N1").Value = " TARGET COMPLETION DATE "
$ws.Range("O1").Value = "ACTUAL DATE OF COMPLETION"
$ws.Range("P1").Value = "PROJECT ID"
$ws.Range("Q1").Value = "CONTRACT ID"
$ws.Range("R1").Value = "ISSUANCE OF INVITATION TO BID"
$ws.Range("S1").Value = "PRE-SUBMISSION CONFERENCE"
$ws.Range("T1").Value = "BID OPENING"
$ws.Range("U1").Value = "ISSUANCE OF RESOLUTION TO AWARD"
$ws.Range("V1").Value = "ISSUANCE OF NOTICE TO PROCEED"
$ws.Range("W1").Value = "NAME OF CONTRACTOR"
$ws.Range("X1").Value = "OTHER REMARKS"
